$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing data row (1084) down
# across the new rows 1085:1121 so column A keeps the date style (s="2")
# and columns B:F remain unstyled, matching the source rows.
$src = $ws.Range("A1084:F1084")
$dst = $ws.Range("A1085:F1121")
$src.Copy()
$dst.PasteSpecial(-4122)

# Populate the new rows with the OHLCV data
$data = @(
    @(1085, 45534.5, 7.752, 7.86, 7.383, 7.415, 531080.21),
    @(1086, 45534.66666666666, 7.419, 7.62, 7.336, 7.537, 353036.67),
    @(1087, 45534.83333333334, 7.533, 7.656, 7.52, 7.617, 112784.17),
    @(1088, 45535, 7.618, 7.679, 7.595, 7.654, 174660.9),
    @(1089, 45535.16666666666, 7.654, 7.722, 7.611, 7.618, 115550.65),
    @(1090, 45535.33333333334, 7.618, 7.631, 7.519, 7.523, 87657.78),
    @(1091, 45535.5, 7.54, 7.558, 7.457, 7.537, 99763.42999999999),
    @(1092, 45535.66666666666, 7.537, 7.562, 7.435, 7.5, 88925.3),
    @(1093, 45535.83333333334, 7.501, 7.53, 7.487, 7.525, 17300.59),
    @(1094, 45536, 7.56, 7.583, 7.41, 7.481, 85597.41),
    @(1095, 45536.16666666666, 7.482, 7.497, 7.351, 7.436, 119813.92),
    @(1096, 45536.33333333334, 7.437, 7.449, 7.352, 7.354, 41074.14),
    @(1097, 45536.5, 7.359, 7.436, 7.189, 7.402, 163731.73),
    @(1098, 45536.66666666666, 7.404, 7.451, 7.245, 7.43, 186237.45),
    @(1099, 45536.83333333334, 7.425, 7.462, 7.087, 7.191, 207139.98),
    @(1100, 45537, 7.184, 7.241, 7.131, 7.224, 158287.26),
    @(1101, 45537.16666666666, 7.225, 7.311, 7.115, 7.146, 157363.35),
    @(1102, 45537.33333333334, 7.145, 7.427, 7.144, 7.36, 301114.78),
    @(1103, 45537.5, 7.363, 7.38, 7.309, 7.315, 23523.98),
    @(1104, 45537.66666666666, 7.391, 7.435, 7.342, 7.387, 141472.7),
    @(1105, 45537.83333333334, 7.393, 7.394, 7.351, 7.388, 18164.48),
    @(1106, 45538, 7.524, 7.645, 7.493, 7.534, 227185.78),
    @(1107, 45538.16666666666, 7.533, 7.555, 7.469, 7.518, 122698.47),
    @(1108, 45538.33333333334, 7.517, 7.543, 7.407, 7.475, 170621.45),
    @(1109, 45538.5, 7.47, 7.519, 7.149, 7.164, 331364.9),
    @(1110, 45538.66666666666, 7.161, 7.269, 7.151, 7.252, 135147.86),
    @(1111, 45538.83333333334, 7.249, 7.293, 7.144, 7.162, 69410.67),
    @(1112, 45539, 7.158, 7.22, 6.9, 7.191, 436059.18),
    @(1113, 45539.16666666666, 7.191, 7.279, 7.16, 7.253, 234856.98),
    @(1114, 45539.33333333334, 7.255, 7.31, 7.151, 7.217, 232717.75),
    @(1115, 45539.5, 7.216, 7.573, 7.171, 7.562, 491769.94),
    @(1116, 45539.66666666666, 7.562, 7.76, 7.522, 7.539, 320665.43),
    @(1117, 45539.83333333334, 7.555, 7.687, 7.549, 7.557, 163991.16),
    @(1118, 45540, 7.557, 7.645, 7.441, 7.466, 207046.71),
    @(1119, 45540.16666666666, 7.466, 7.499, 7.362, 7.436, 393798.24),
    @(1120, 45540.33333333334, 7.435, 7.464, 7.328, 7.333, 376957.76),
    @(1121, 45540.5, 7.333, 7.335, 7.33, 7.335, 972.08)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}

Write-Host "Added rows 1085-1121"